$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.861.57'
$ws.Range("E2").Value = '  -0.80%  '
$ws.Range("D3").Value = '2.218.03'
$ws.Range("E3").Value = '  -1.18%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '256.54'
$ws.Range("E5").Value = '  +2.37%  '
$ws.Range("B6").Value = 'XRP'
$ws.Range("C6").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D6").Value = '0.618'
$ws.Range("E6").Value = '  +0.27%  '
$ws.Range("B7").Value = 'Solana'
$ws.Range("C7").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D7").Value = '77.79'
$ws.Range("E7").Value = '  +3.27%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("E9").Value = '  -1.21%  '
$ws.Range("D10").Value = '43.08'
$ws.Range("E10").Value = '  +4.52%  '
$ws.Range("D11").Value = '0.0910'
$ws.Range("E11").Value = '  -2.54%  '
$ws.Range("E12").Value = '  +0.64%  '
$ws.Range("D13").Value = '0.102'
$ws.Range("E13").Value = '  +0.52%  '
$ws.Range("D14").Value = '2.551.79'
$ws.Range("E14").Value = '  -1.14%  '
$ws.Range("D15").Value = '14.46'
$ws.Range("E15").Value = '  -1.50%  '
$ws.Range("D16").Value = '2.218.85'
$ws.Range("E16").Value = '  -1.20%  '
$ws.Range("D17").Value = '0.783'
$ws.Range("E17").Value = '  -1.31%  '
$ws.Range("D18").Value = '42.823.20'
$ws.Range("E18").Value = '  -0.55%  '
$ws.Range("E19").Value = '  -1.73%  '
$ws.Range("E20").Value = '  -0.27%  '
$ws.Range("D21").Value = '5.98'
$ws.Range("E21").Value = '  -0.56%  '
$ws.Range("D22").Value = '2.28'
$ws.Range("E22").Value = '  +2.97%  '
$ws.Range("D23").Value = '229.89'
$ws.Range("E23").Value = '  -0.04%  '
$ws.Range("D24").Value = '9.30'
$ws.Range("E24").Value = '  -4.82%  '
$ws.Range("E25").Value = '  -0.07%  '
$ws.Range("D26").Value = '42.70'
$ws.Range("E26").Value = '  +7.66%  '
$ws.Range("D27").Value = '10.76'
$ws.Range("E27").Value = '  -0.50%  '
$ws.Range("E28").Value = '  -2.71%  '
$ws.Range("D29").Value = '2.22'
$ws.Range("E29").Value = '  -0.11%  '
$ws.Range("D30").Value = '2.20'
$ws.Range("E30").Value = '  -2.75%  '
$ws.Range("D31").Value = '174.04'
$ws.Range("E31").Value = '  +1.23%  '
$ws.Range("E32").Value = '  +0.82%  '
$ws.Range("D33").Value = '0.0871'
$ws.Range("E33").Value = '  +8.33%  '
$ws.Range("D34").Value = '5.21'
$ws.Range("E34").Value = '  -1.04%  '
$ws.Range("E35").Value = '  -0.36%  '
$ws.Range("D36").Value = '0.0357'
$ws.Range("E36").Value = '  +7.75%  '
$ws.Range("E37").Value = '  -2.27%  '
$ws.Range("D38").Value = '4.32'
$ws.Range("E38").Value = '  -3.51%  '
$ws.Range("D39").Value = '13.01'
$ws.Range("E39").Value = '  -1.34%  '
$ws.Range("D40").Value = '2.82'
$ws.Range("E40").Value = '  +16.30%  '
$ws.Range("E41").Value = '  -0.21%  '
$ws.Range("D42").Value = '61.20'
$ws.Range("E42").Value = '  +2.90%  '
$ws.Range("E43").Value = '  -2.67%  '
$ws.Range("D44").Value = '5.32'
$ws.Range("E44").Value = '  -2.37%  '
$ws.Range("E45").Value = '  +0.85%  '
$ws.Range("D46").Value = '103.72'
$ws.Range("E46").Value = '  -0.76%  '
$ws.Range("D47").Value = '8.43'
$ws.Range("E47").Value = '  -3.01%  '
$ws.Range("D48").Value = '0.0971'
$ws.Range("E48").Value = '  -2.38%  '
$ws.Range("D49").Value = '1.12'
$ws.Range("E50").Value = '  -1.51%  '
$ws.Range("D51").Value = '1.46'
$ws.Range("E51").Value = '  +20.03%  '
